# Append the June 2021 liquidity override row (19437/60001 June 2021 report)
# to the bottom of the existing table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data row: 2021-06-30, YMM US Equity, L0
$ws.Range("A6").Value = 44377          # serial date for 2021-06-30
$ws.Range("B6").Value = "YMM US Equity"
$ws.Range("C6").Value = "L0"

# Copy formatting (date number format, etc.) from the row above so the new
# row matches the existing style instead of Excel auto-creating a new one.
$ws.Range("A5:C5").Copy()
$ws.Range("A6:C6").PasteSpecial(-4122)  # xlPasteFormats
